$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.189129442330756
$ws.Range("C2").Value = 4.975048808213939
$ws.Range("D2").Value = 5.512714576635176
$ws.Range("F2").Value = 29.70730693473807
$ws.Range("G2").Value = 38.57026274743739
$ws.Range("H2").Value = 16.36177310450883
$ws.Range("K2").Value = 7.565589750431634
$ws.Range("M2").Value = 19.97959944550328
$ws.Range("B3").Value = 8.009700938023469
$ws.Range("C3").Value = 4.906111319743714
$ws.Range("D3").Value = 5.504093905534972
$ws.Range("F3").Value = 29.51843853347945
$ws.Range("G3").Value = 38.19099813884709
$ws.Range("H3").Value = 16.35037350796008
$ws.Range("K3").Value = 7.471119177851576
$ws.Range("M3").Value = 19.39227238399427
$ws.Range("B4").Value = 7.900105355938185
$ws.Range("C4").Value = 4.862197406099911
$ws.Range("D4").Value = 5.498386880100088
$ws.Range("F4").Value = 29.40835852569619
$ws.Range("G4").Value = 37.96536710553634
$ws.Range("H4").Value = 16.34585612347178
$ws.Range("K4").Value = 7.414807240002277
$ws.Range("M4").Value = 19.03082703143031
$ws.Range("B5").Value = 7.855665566999488
$ws.Range("C5").Value = 4.843910346781212
$ws.Range("D5").Value = 5.495957399138905
$ws.Range("F5").Value = 29.36501117472244
$ws.Range("G5").Value = 37.8753220063506
$ws.Range("H5").Value = 16.34463963374358
$ws.Range("K5").Value = 7.392317565478507
$ws.Range("M5").Value = 18.88357286818615
$ws.Range("B6").Value = 7.848302016877367
$ws.Range("C6").Value = 4.840850343404792
$ws.Range("D6").Value = 5.49554772260739
$ws.Range("F6").Value = 29.35790550623998
$ws.Range("G6").Value = 37.86048720968738
$ws.Range("H6").Value = 16.34447534918797
$ws.Range("K6").Value = 7.388611748882731
$ws.Range("M6").Value = 18.8591310763805
$ws.Range("B7").Value = 7.899505026602559
$ws.Range("C7").Value = 4.861952357230463
$ws.Range("D7").Value = 5.498354535303744
$ws.Range("F7").Value = 29.40776777043219
$ws.Range("G7").Value = 37.96414492322884
$ws.Range("H7").Value = 16.34583718923477
$ws.Range("K7").Value = 7.414502040597768
$ws.Range("M7").Value = 19.0288406256858
$ws.Range("B8").Value = 8.127188920104448
$ws.Range("C8").Value = 4.951612458056958
$ws.Range("D8").Value = 5.509828226910803
$ws.Range("F8").Value = 29.64098289819284
$ws.Range("G8").Value = 38.43803976865338
$ws.Range("H8").Value = 16.35732761968096
$ws.Range("K8").Value = 7.532684330298892
$ws.Range("M8").Value = 19.77740357007992
$ws.Range("B9").Value = 8.574952593234181
$ws.Range("C9").Value = 5.114560196606695
$ws.Range("D9").Value = 5.529032744023445
$ws.Range("F9").Value = 30.14353813622148
$ws.Range("G9").Value = 39.42103391901561
$ws.Range("H9").Value = 16.39953281665012
$ws.Range("K9").Value = 7.776430719058661
$ws.Range("M9").Value = 21.22865725170483
$ws.Range("B10").Value = 8.900498567885789
$ws.Range("C10").Value = 5.226064144548062
$ws.Range("D10").Value = 5.541118178412983
$ws.Range("F10").Value = 30.53819639073357
$ws.Range("G10").Value = 40.17072680969137
$ws.Range("H10").Value = 16.44249111273435
$ws.Range("K10").Value = 7.960857675782108
$ws.Range("M10").Value = 22.27217420542573
$ws.Range("B11").Value = 9.047031351841518
$ws.Range("C11").Value = 5.27493598628023
$ws.Range("D11").Value = 5.546173376439569
$ws.Range("F11").Value = 30.72276173942756
$ws.Range("G11").Value = 40.51658717414697
$ws.Range("H11").Value = 16.46461339553149
$ws.Range("K11").Value = 8.045501854050201
$ws.Range("M11").Value = 22.73962633856685
$ws.Range("B12").Value = 9.10223206861477
$ws.Range("C12").Value = 5.29317061917857
$ws.Range("D12").Value = 5.548023734245501
$ws.Range("F12").Value = 30.79333116781981
$ws.Range("G12").Value = 40.64814918442116
$ws.Range("H12").Value = 16.47335964071033
$ws.Range("K12").Value = 8.077627553761022
$ws.Range("M12").Value = 22.91542357959502
$ws.Range("B13").Value = 9.090357388909267
$ws.Range("C13").Value = 5.289255653747897
$ws.Range("D13").Value = 5.547628074957429
$ws.Range("F13").Value = 30.7781033351358
$ws.Range("G13").Value = 40.61979022671118
$ws.Range("H13").Value = 16.47145960398304
$ws.Range("K13").Value = 8.070706017081921
$ws.Range("M13").Value = 22.87761921164687
$ws.Range("B14").Value = 9.05157887160213
$ws.Range("C14").Value = 5.276441645181461
$ws.Range("D14").Value = 5.546326889596108
$ws.Range("F14").Value = 30.72855419379587
$ws.Range("G14").Value = 40.52739957176762
$ws.Range("H14").Value = 16.4653255715369
$ws.Range("K14").Value = 8.048143592136732
$ws.Range("M14").Value = 22.75411460787089
$ws.Range("B15").Value = 9.027786574420926
$ws.Range("C15").Value = 5.268557097123871
$ws.Range("D15").Value = 5.545521540666606
$ws.Range("F15").Value = 30.69829091370179
$ws.Range("G15").Value = 40.47088174228934
$ws.Range("H15").Value = 16.46161629056621
$ws.Range("K15").Value = 8.034331916688142
$ws.Range("M15").Value = 22.6783011277896
$ws.Range("B16").Value = 8.890885371662197
$ws.Range("C16").Value = 5.222832489211533
$ws.Range("D16").Value = 5.540778878853034
$ws.Range("F16").Value = 30.5262320578689
$ws.Range("G16").Value = 40.14821240682629
$ws.Range("H16").Value = 16.44109712616659
$ws.Range("K16").Value = 7.95533803881178
$ws.Range("M16").Value = 22.2414639559497
$ws.Range("B17").Value = 8.806454721685125
$ws.Range("C17").Value = 5.194303402292612
$ws.Range("D17").Value = 5.537755853909266
$ws.Range("F17").Value = 30.42193692999088
$ws.Range("G17").Value = 39.95142831354568
$ws.Range("H17").Value = 16.42916882667496
$ws.Range("K17").Value = 7.907044397732228
$ws.Range("M17").Value = 21.9714925015245
$ws.Range("B18").Value = 8.757749189576417
$ws.Range("C18").Value = 5.177720244449055
$ws.Range("D18").Value = 5.535975455289247
$ws.Range("F18").Value = 30.3624256961909
$ws.Range("G18").Value = 39.83870321519174
$ws.Range("H18").Value = 16.4225510003817
$ws.Range("K18").Value = 7.879339564863049
$ws.Range("M18").Value = 21.81553910994276
$ws.Range("B19").Value = 8.741235741541439
$ws.Range("C19").Value = 5.172075772021401
$ws.Range("D19").Value = 5.535365501261222
$ws.Range("F19").Value = 30.34235942137827
$ws.Range("G19").Value = 39.80061849534605
$ws.Range("H19").Value = 16.42035211293004
$ws.Range("K19").Value = 7.869972685837692
$ws.Range("M19").Value = 21.76262597787839
$ws.Range("B20").Value = 8.815457810366974
$ws.Range("C20").Value = 5.197358421616639
$ws.Range("D20").Value = 5.538081971278203
$ws.Range("F20").Value = 30.4329903251696
$ws.Range("G20").Value = 39.97232944851805
$ws.Range("H20").Value = 16.43041347698141
$ws.Range("K20").Value = 7.912178099442316
$ws.Range("M20").Value = 22.00030238771085
$ws.Range("B21").Value = 9.06297736750448
$ws.Range("C21").Value = 5.280212858515424
$ws.Range("D21").Value = 5.546710816888047
$ws.Range("F21").Value = 30.74308992959408
$ws.Range("G21").Value = 40.55452168212555
$ws.Range("H21").Value = 16.46711728834562
$ws.Range("K21").Value = 8.054769020816671
$ws.Range("M21").Value = 22.7904252323536
$ws.Range("B22").Value = 9.223040055280027
$ws.Range("C22").Value = 5.332774723016742
$ws.Range("D22").Value = 5.551977457390414
$ws.Range("F22").Value = 30.94969110867336
$ws.Range("G22").Value = 40.93842272774839
$ws.Range("H22").Value = 16.49325509288551
$ws.Range("K22").Value = 8.148369137949871
$ws.Range("M22").Value = 23.29965495391865
$ws.Range("B23").Value = 9.137787624001895
$ws.Range("C23").Value = 5.304868620718935
$ws.Range("D23").Value = 5.54920076424826
$ws.Range("F23").Value = 30.83907956528379
$ws.Range("G23").Value = 40.73324943922639
$ws.Range("H23").Value = 16.47910890795011
$ws.Range("K23").Value = 8.098386865236956
$ws.Range("M23").Value = 23.02857802105428
$ws.Range("B24").Value = 8.811388023732739
$ws.Range("C24").Value = 5.195977811977459
$ws.Range("D24").Value = 5.537934665725751
$ws.Range("F24").Value = 30.42799168193747
$ws.Range("G24").Value = 39.96287876473131
$ws.Range("H24").Value = 16.42985002315545
$ws.Range("K24").Value = 7.909856963919461
$ws.Range("M24").Value = 21.9872797396756
$ws.Range("B25").Value = 8.454130528366205
$ws.Range("C25").Value = 5.071895815273809
$ws.Range("D25").Value = 5.524191773577717
$ws.Range("F25").Value = 30.00295455888697
$ws.Range("G25").Value = 39.14989116380256
$ws.Range("H25").Value = 16.38601225687796
$ws.Range("K25").Value = 7.709409397610441
$ws.Range("M25").Value = 20.83919594976571
